$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has blank placeholder rows interleaved with the
# real header/data rows (row 1 blank, row 2 blank, header on row 3,
# blank row 4, data on rows 5-17, blank row 18, blank footer row 19).
# Removing the blank spacer rows collapses everything upward so the
# header lands on row 1, the data on rows 2-14, and the footer on row 16.
#
# Delete from the bottom up (by original row number) so earlier deletes
# don't invalidate the row numbers used by later deletes.
$ws.Rows(4).Delete()
$ws.Rows(2).Delete()
$ws.Rows(1).Delete()
